# Update column G ("K") values on Sheet1 to reflect the regenerated
# save_data (switch from Strike# to K, recalculated std/mean, s_vals).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$newValues = @{
    2  = 1
    3  = 0
    4  = 0
    5  = 2
    6  = 0
    7  = 1
    8  = 1
    9  = 2
    10 = 1
    11 = 2
    12 = 1
    13 = 1
    14 = 0
    15 = 3
    16 = 2
    17 = 0
    18 = 0
    19 = 3
    20 = 0
    21 = 2
    22 = 3
    23 = 1
    25 = 1
}

foreach ($row in $newValues.Keys) {
    $ws.Range("G$row").Value = $newValues[$row]
}
